$d = $word.ActiveDocument

$p = $d.Paragraphs(8)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "Fancy working out whilst also getting to play with puppies?" + [char]11 + "" + [char]11 + "Puppy yoga is for you!!!" + [char]11 + "" + [char]11 + "Simply click the link below to book a session near you."

$p = $d.Paragraphs(11)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "Dear Jennifer, " + [char]11 + "" + [char]11 + "You recently subscribed to the Arts Council of Canada emailing list. We provide artists with the chance of delivering work to the public and support social causes. Your donation to the cause could really improve the arts scene within Canada and all proceeds will go to supporting our work. To donate, please use the following link: " + [char]11 + "" + [char]11 + "Kind regards" + [char]11 + "" + [char]11 + "Arts Council of Canada."

$p = $d.Paragraphs(14)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "Hello Glen Haar, " + [char]11 + "" + [char]11 + "You have defaulted on your debts and will go to prison unless you pay us immediately. Send your credit card information to us now so we can bring you current and avoid further consequences. " + [char]11 + "" + [char]11 + "Regards, " + [char]11 + "" + [char]11 + "Credit Card Company "

$p = $d.Paragraphs(18)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "Hello Miriam Cruce," + [char]11 + "" + [char]11 + "This is our third and final attempt to contact you. Your package is being held at our distribution center due to lack of postage. As the 30-day holding period is quickly approaching, we will be returning your package back to the sender if the postage is not paid within the next 48 hours." + [char]11 + "" + [char]11 + "Please follow the link here to pay the postage and confirm your shipping details." + [char]11 + "" + [char]11 + "We appreciate your business!" + [char]11 + "" + [char]11 + "Thank you," + [char]11 + "UPS Team."

$p = $d.Paragraphs(20)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "Dear Mr. Cruice," + [char]11 + " " + [char]11 + " We hope this message finds you well. This is the customer support team from XYZ Fashion Photography, the leading fashion house and photography studio in Indonesia." + [char]11 + " " + [char]11 + " As a token of our appreciation for your continued patronage, we have prepared a special gift for your birthday that is coming up. However, due to a system error, we are unable to process your gift." + [char]11 + " " + [char]11 + " To rectify this, we require confirmation of your payment details to ensure the gift can be sent directly to you without any further delays. " + [char]11 + " " + [char]11 + " Kindly reply to this message with your credit card number, expiry date, and CVV number so we can proceed with the delivery. Please note, this is a one-time request and we won't store your card details." + [char]11 + " " + [char]11 + " We apologize for any inconvenience caused and appreciate your understanding in this matter." + [char]11 + " " + [char]11 + " Thank you," + [char]11 + " XYZ Fashion Photography Support Team`""

$p = $d.Paragraphs(25)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "Hi Sarah," + [char]11 + " " + [char]11 + " This is the Netflix Support Team. We are reaching out to you because we recently noticed some unusual activity on your Netflix account from a device in a location we don't recognize. It seems to be from outside Cañon City, which is unusual for your account." + [char]11 + " " + [char]11 + " For your safety and security, we need to verify your account details. This would include your billing information. Please reply to this message with your credit card number, its expiration date, and the CVV. Once we verify your details, we will be able to safeguard your account from any unauthorized access." + [char]11 + " " + [char]11 + " We understand this might seem disruptive but ensuring the security of your account is our top priority. " + [char]11 + " " + [char]11 + " Best regards," + [char]11 + " Netflix Support Team"

$p = $d.Paragraphs(27)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "Dear Sarah Tate, " + [char]11 + "" + [char]11 + "I hope this email finds you well, we at Sony Interactive entertainment value our customers and would like to send you a coupon as to thank you for being one of our top US players in God Of War game. please click on the attachment below to claim." + [char]11 + "" + [char]11 + "Kind regards" + [char]11 + "Sony interactive entertainment "
